$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in rows 9-14 with Agent / OCM Agent Summary Report / Date Range (columns A-C)
# and the refresh-interval numeric values in column I (stored as text, quote-prefixed
# like the existing I8 cell), matching the pattern already present in row 8.
$iValues = @{ 9 = "30"; 10 = "2"; 11 = "23"; 12 = "2"; 13 = "14"; 14 = "2" }

for ($r = 9; $r -le 14; $r++) {
    $ws.Range("A$r").NumberFormat = "@"
    $ws.Range("A$r").Value = "Agent"

    $ws.Range("B$r").NumberFormat = "@"
    $ws.Range("B$r").Value = "OCM Agent Summary Report"

    $ws.Range("C$r").Value = "Date Range"

    $ws.Range("I$r").Value = "'" + $iValues[$r]
}

# Move the active selection from D2 to H15, as recorded in the saved workbook view.
$ws.Range("H15").Select()
